{"js": "// Apply the \"Incas Return\" review copy edits described in the diff.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Helper: replace the text of the paragraph at `index`, after verifying\n// that it currently holds `expected` text (keeps the script honest if the\n// document shape ever drifts from what we expect).\nfunction replaceAt(index, expected, replacement) {\n  const p = items[index];\n  if (p.text !== expected) {\n    throw new Error(\n      \"Paragraph \" + index + \" text mismatch. Expected: '\" + expected + \"' Actual: '\" + p.text + \"'\"\n    );\n  }\n  p.insertText(replacement, \"Replace\");\n  return p;\n}\n\n// 1) Title heading (Heading1, paragraph 0).\nreplaceAt(0, \"Play Incas Return Slot for Free - Game Review\", \"Play Incas Return Free | Online Slot Review\");\n\n// 2) \"What we like\" bullet list (paragraphs 35-38).\nreplaceAt(35, \"Wild symbol replaces any other symbol\", \"Standard slot game following online slot standards\");\n\nconst scatterBullet = replaceAt(\n  36,\n  \"Scatter symbol grants access to Bonus Game\",\n  \"Wild symbol and Scatter symbol with Bonus Game\"\n);\n\n// New bullet inserted right after the (now renamed) Wild/Scatter bullet.\nscatterBullet.insertParagraph(\"Pre-Columbian civilization theme\", \"After\");\n\n// \"Illustrative symbols rich in details\" (paragraph 37) is unchanged.\n\n// Old \"Standard slot follows online slot standards set by Cristaltec\" bullet\n// (paragraph 38) is removed entirely.\nif (items[38].text !== \"Standard slot follows online slot standards set by Cristaltec\") {\n  throw new Error(\"Paragraph 38 text mismatch. Actual: '\" + items[38].text + \"'\");\n}\nitems[38].delete();\n\n// 3) \"What we don't like\" bullet list (paragraphs 40-41).\nreplaceAt(40, \"Weak music loop quality\", \"Poorly looped music\");\nreplaceAt(41, \"No sound effects\", \"Difficulty distinguishing symbols\");\n\n// 4) Bold summary title near the end of the document (paragraph 42).\nreplaceAt(\n  42,\n  \"Play Incas Return Slot for Free - Game Review\",\n  \"Play Incas Return Free | Online Slot Review\"\n);\n\n// 5) Italic summary/description paragraph (paragraph 43).\nreplaceAt(\n  43,\n  \"Read our unbiased review of the Incas Return slot game and play for free. Discover its gameplay, theme, technical features, and audio components.\",\n  \"Read our review of Incas Return and play this slot game for free.\"\n);\n\nawait context.sync();\n", "ps1": "# Apply the \"Incas Return\" review copy edits described in the diff.\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($index) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (CR); strip it\n    # for clean comparisons.\n    return $d.Paragraphs.Item($index).Range.Text.TrimEnd([char]13)\n}\n\nfunction Set-ParaText($index, $expected, $replacement) {\n    $actual = Get-ParaText $index\n    if ($actual -ne $expected) {\n        throw \"Paragraph $index text mismatch. Expected: '$expected' Actual: '$actual'\"\n    }\n    $r = $d.Paragraphs.Item($index).Range\n    $r.MoveEnd(1, -1)  # wdCharacter; exclude the trailing paragraph mark\n    $r.Text = $replacement\n}\n\n# 1) Title heading (Heading 1, paragraph 1).\nSet-ParaText 1 \"Play Incas Return Slot for Free - Game Review\" \"Play Incas Return Free | Online Slot Review\"\n\n# 2) \"What we like\" bullet list (paragraphs 36-39).\nSet-ParaText 36 \"Wild symbol replaces any other symbol\" \"Standard slot game following online slot standards\"\nSet-ParaText 37 \"Scatter symbol grants access to Bonus Game\" \"Wild symbol and Scatter symbol with Bonus Game\"\n\n# New bullet inserted right after the (now renamed) Wild/Scatter bullet.\n$scatterRange = $d.Paragraphs.Item(37).Range\n$scatterRange.InsertParagraphAfter()\n$newBullet = $d.Paragraphs.Item(38).Range\n$newBullet.Text = \"Pre-Columbian civilization theme\"\n\n# \"Illustrative symbols rich in details\" (now paragraph 39) is unchanged.\n\n# Old \"Standard slot follows online slot standards set by Cristaltec\" bullet\n# (now paragraph 40) is removed entirely.\n$cristaltecText = Get-ParaText 40\nif ($cristaltecText -ne \"Standard slot follows online slot standards set by Cristaltec\") {\n    throw \"Paragraph 40 text mismatch. Actual: '$cristaltecText'\"\n}\n$d.Paragraphs.Item(40).Range.Delete()\n\n# 3) \"What we don't like\" bullet list (paragraphs 41-42 after the insert/delete above).\nSet-ParaText 41 \"Weak music loop quality\" \"Poorly looped music\"\nSet-ParaText 42 \"No sound effects\" \"Difficulty distinguishing symbols\"\n\n# 4) Bold summary title near the end of the document (paragraph 43).\nSet-ParaText 43 \"Play Incas Return Slot for Free - Game Review\" \"Play Incas Return Free | Online Slot Review\"\n\n# 5) Italic summary/description paragraph (paragraph 44).\nSet-ParaText 44 \"Read our unbiased review of the Incas Return slot game and play for free. Discover its gameplay, theme, technical features, and audio components.\" \"Read our review of Incas Return and play this slot game for free.\"\n"}
